$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(137, 1).Value = "espino"
$ws.Cells.Item(137, 2).Value = "paolo"
$ws.Cells.Item(137, 3).Value = 502179
$ws.Cells.Item(137, 4).Value = "espip001"
$ws.Cells.Item(137, 5).Value = "espinpa01"
$ws.Cells.Item(137, 6).Value = 8246
$ws.Cells.Item(137, 7).Value = 2017
$ws.Cells.Item(137, 8).Value = 2021
$ws.Cells.Item(137, 9).Value = 29989

$ws.Cells.Item(138, 1).Value = "stroman"
$ws.Cells.Item(138, 2).Value = "marcus"
$ws.Cells.Item(138, 3).Value = 573186
$ws.Cells.Item(138, 4).Value = "strom001"
$ws.Cells.Item(138, 5).Value = "stromma01"
$ws.Cells.Item(138, 6).Value = 13431
$ws.Cells.Item(138, 7).Value = 2014
$ws.Cells.Item(138, 8).Value = 2021
$ws.Cells.Item(138, 9).Value = 32815

$ws.Cells.Item(139, 1).Value = "pineda"
$ws.Cells.Item(139, 2).Value = "michael"
$ws.Cells.Item(139, 3).Value = 501381
$ws.Cells.Item(139, 4).Value = "pinem001"
$ws.Cells.Item(139, 5).Value = "pinedmi01"
$ws.Cells.Item(139, 6).Value = 5372
$ws.Cells.Item(139, 7).Value = 2011
$ws.Cells.Item(139, 8).Value = 2021
$ws.Cells.Item(139, 9).Value = 30937

$ws.Cells.Item(140, 1).Value = "foltynewicz"
$ws.Cells.Item(140, 2).Value = "mike"
$ws.Cells.Item(140, 3).Value = 592314
$ws.Cells.Item(140, 4).Value = "foltm001"
$ws.Cells.Item(140, 5).Value = "foltymi01"
$ws.Cells.Item(140, 6).Value = 10811
$ws.Cells.Item(140, 7).Value = 2014
$ws.Cells.Item(140, 8).Value = 2021
$ws.Cells.Item(140, 9).Value = 31819

$ws.Cells.Item(141, 1).Value = "garcia"
$ws.Cells.Item(141, 2).Value = "luis"
$ws.Cells.Item(141, 3).Value = 671277
$ws.Cells.Item(141, 4).Value = "garcl006"
$ws.Cells.Item(141, 5).Value = "garcilu04"
$ws.Cells.Item(141, 6).Value = -1
$ws.Cells.Item(141, 7).Value = 2020
$ws.Cells.Item(141, 8).Value = 2021
$ws.Cells.Item(141, 9).Value = 4684365

$ws.Cells.Item(142, 1).Value = "ohtani"
$ws.Cells.Item(142, 2).Value = "shohei"
$ws.Cells.Item(142, 3).Value = 660271
$ws.Cells.Item(142, 4).Value = "ohtas001"
$ws.Cells.Item(142, 5).Value = "ohtansh01"
$ws.Cells.Item(142, 6).Value = 19755
$ws.Cells.Item(142, 7).Value = 2018
$ws.Cells.Item(142, 8).Value = 2021
$ws.Cells.Item(142, 9).Value = 39832

$ws.Cells.Item(143, 1).Value = "smith"
$ws.Cells.Item(143, 2).Value = "caleb"
$ws.Cells.Item(143, 3).Value = 592761
$ws.Cells.Item(143, 4).Value = "smitc006"
$ws.Cells.Item(143, 5).Value = "smithca03"
$ws.Cells.Item(143, 6).Value = 14875
$ws.Cells.Item(143, 7).Value = 2017
$ws.Cells.Item(143, 8).Value = 2021
$ws.Cells.Item(143, 9).Value = 36081

$ws.Cells.Item(144, 1).Value = "pivetta"
$ws.Cells.Item(144, 2).Value = "nick"
$ws.Cells.Item(144, 3).Value = 601713
$ws.Cells.Item(144, 4).Value = "piven001"
$ws.Cells.Item(144, 5).Value = "pivetni01"
$ws.Cells.Item(144, 6).Value = 15454
$ws.Cells.Item(144, 7).Value = 2017
$ws.Cells.Item(144, 8).Value = 2021
$ws.Cells.Item(144, 9).Value = 36071
